$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the value out of B6 (previously held the shared string "ffffff")
$ws.Range("B6").ClearContents()

# Move the active selection from B6 to E8
$ws.Range("E8").Select()
